$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column C
$ws.Range("C1").Value = "Similar"

# Column C "Similar" annotations for the relevant rows
$ws.Range("C2").Value = "2, 5"
$ws.Range("C4").Value = "4, 14, 28, 44"
$ws.Range("C5").Value = "5, 29"
$ws.Range("C6").Value = "6, 27, 37"
$ws.Range("C8").Value = 15
$ws.Range("C14").Value = "4, 14, 28, 44"
$ws.Range("C25").Value = "2, 25"
$ws.Range("C27").Value = "6, 27, 37"

# Row 29 gains a Decision value and a Similar annotation
$ws.Range("B29").Value = "yes"
$ws.Range("C29").Value = "5, 29"

$ws.Range("C44").Value = "4, 14, 28, 44"

# Update print/page setup (portrait orientation)
$ws.PageSetup.Orientation = 1

# Update the active cell selection to match the authored state
$ws.Range("D14").Select()
